$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.3796008642440598
$ws.Range("C2").Value = 3.6918003194439302
$ws.Range("E2").Value = 1.5780000000000001

$ws.Range("B3").Value = 1.55897105099264
$ws.Range("C3").Value = 1.1699381174883601
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1.2589999999999999

$ws.Range("B4").Value = 3.5859652166517799
$ws.Range("C4").Value = 3.6116202332824399
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1.389

$ws.Range("B5").Value = 1.9287925424064001
$ws.Range("C5").Value = 1.7164829247606299
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1.5229999999999999

$ws.Range("B6").Value = 2.02502948617445
$ws.Range("C6").Value = 1.6086746094601101
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1.5069999999999999
$ws.Range("F6").Value = 14.856

$ws.Rows("7:25").Select()
